$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Periodo Mora" column (E) for the two existing debt rows: swap the two
# period values so the first row now shows period 1601 and the second 1602.
$ws.Range("E16").Value = "1601"
$ws.Range("E17").Value = "1602"

# "Salario Basico" column (G): update the worker's base salary for both
# mora periods (was 0, now populated with the real salary amount).
$ws.Range("G16").Value = 644350
$ws.Range("G17").Value = 644350
